# ---------------------------------------------------------------
# Applies the commit "Fixed update to excel issue":
#  1. Rename header B1 on "Weekly Quantity"  -> "Weekly_PO_Qty"
#  2. Rename header B1 on "Monthly Trend"    -> "Monthly_PO_Qty"
#  3. Add a new "PO Forecast" worksheet (sheetId 3 / rId3, placed
#     after "Monthly Trend") with forecast data: ds, PO_Forecast,
#     yhat_lower, yhat_upper
# ---------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1 & 2: rename the "Requested quantity" headers ---------------
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("B1").Value = "Weekly_PO_Qty"

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: add the new "PO Forecast" sheet, positioned after the -----
#        last existing sheet so it lands as sheetId 3 / rId3 -------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "PO Forecast"

# Header row
$ws3.Cells.Item(1,1).Value = "ds"
$ws3.Cells.Item(1,2).Value = "PO_Forecast"
$ws3.Cells.Item(1,3).Value = "yhat_lower"
$ws3.Cells.Item(1,4).Value = "yhat_upper"

# Data rows (2-18)
$ws3.Cells.Item(2,1).Value = 44948.99999999999
$ws3.Cells.Item(2,2).Value = 14
$ws3.Cells.Item(2,3).Value = 4.811050398242077
$ws3.Cells.Item(2,4).Value = 21.45061053583695
$ws3.Cells.Item(3,1).Value = 44955.99999999999
$ws3.Cells.Item(3,2).Value = 14
$ws3.Cells.Item(3,3).Value = 5.287684375463265
$ws3.Cells.Item(3,4).Value = 22.61259161860073
$ws3.Cells.Item(4,1).Value = 44962.99999999999
$ws3.Cells.Item(4,2).Value = 14
$ws3.Cells.Item(4,3).Value = 4.429021688167442
$ws3.Cells.Item(4,4).Value = 22.00203184940514
$ws3.Cells.Item(5,1).Value = 44976.99999999999
$ws3.Cells.Item(5,2).Value = 14
$ws3.Cells.Item(5,3).Value = 5.325844373632378
$ws3.Cells.Item(5,4).Value = 22.3504726283629
$ws3.Cells.Item(6,1).Value = 44990.99999999999
$ws3.Cells.Item(6,2).Value = 14
$ws3.Cells.Item(6,3).Value = 5.433036685267878
$ws3.Cells.Item(6,4).Value = 22.3878063519464
$ws3.Cells.Item(7,1).Value = 45032.99999999999
$ws3.Cells.Item(7,2).Value = 15
$ws3.Cells.Item(7,3).Value = 5.893652005807337
$ws3.Cells.Item(7,4).Value = 23.8695005452282
$ws3.Cells.Item(8,1).Value = 45067.99999999999
$ws3.Cells.Item(8,2).Value = 15
$ws3.Cells.Item(8,3).Value = 6.120525774222306
$ws3.Cells.Item(8,4).Value = 23.99457057624612
$ws3.Cells.Item(9,1).Value = 45081.99999999999
$ws3.Cells.Item(9,2).Value = 15
$ws3.Cells.Item(9,3).Value = 6.819656779480601
$ws3.Cells.Item(9,4).Value = 23.45276419782977
$ws3.Cells.Item(10,1).Value = 45095.99999999999
$ws3.Cells.Item(10,2).Value = 15
$ws3.Cells.Item(10,3).Value = 6.588110681681489
$ws3.Cells.Item(10,4).Value = 23.58978841790612
$ws3.Cells.Item(11,1).Value = 45102.99999999999
$ws3.Cells.Item(11,2).Value = 16
$ws3.Cells.Item(11,3).Value = 7.150655798430816
$ws3.Cells.Item(11,4).Value = 24.70120333675595
$ws3.Cells.Item(12,1).Value = 45109.99999999999
$ws3.Cells.Item(12,2).Value = 16
$ws3.Cells.Item(12,3).Value = 7.28378070461506
$ws3.Cells.Item(12,4).Value = 23.41846609423867
$ws3.Cells.Item(13,1).Value = 45116.99999999999
$ws3.Cells.Item(13,2).Value = 16
$ws3.Cells.Item(13,3).Value = 7.318743449441446
$ws3.Cells.Item(13,4).Value = 23.77995243493174
$ws3.Cells.Item(14,1).Value = 45123.99999999999
$ws3.Cells.Item(14,2).Value = 16
$ws3.Cells.Item(14,3).Value = 6.614419168197328
$ws3.Cells.Item(14,4).Value = 24.03877720827884
$ws3.Cells.Item(15,1).Value = 45130.99999999999
$ws3.Cells.Item(15,2).Value = 16
$ws3.Cells.Item(15,3).Value = 6.608262848741352
$ws3.Cells.Item(15,4).Value = 24.20104640293319
$ws3.Cells.Item(16,1).Value = 45137.99999999999
$ws3.Cells.Item(16,2).Value = 16
$ws3.Cells.Item(16,3).Value = 7.513086684116376
$ws3.Cells.Item(16,4).Value = 24.73359118504321
$ws3.Cells.Item(17,1).Value = 45144.99999999999
$ws3.Cells.Item(17,2).Value = 16
$ws3.Cells.Item(17,3).Value = 7.111461435570476
$ws3.Cells.Item(17,4).Value = 24.63001241059244
$ws3.Cells.Item(18,1).Value = 45151.99999999999
$ws3.Cells.Item(18,2).Value = 16
$ws3.Cells.Item(18,3).Value = 7.701885359177983
$ws3.Cells.Item(18,4).Value = 24.84308134129092

# --- Formatting: reuse the same header / date styles already ------
#     used on the other two sheets (copy formats only, no values) --
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws3.Range("A2:A18").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Restore the originally active sheet (adding a sheet makes it active)
$ws1.Activate()
